# Daily attendance processing - 2025-10-20 15:43:39
# Rotates the "Recorded By" (column G) list of names/emails for each data
# row: the first entry in the comma-separated list is moved to the end.
# Cells that contain only a single value are left unchanged (rotating a
# single-element list is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -isnot [string]) { continue }
    if ($val -eq "") { continue }
    if ($val -eq "Recorded By") { continue }

    $parts = $val -split ", "

    if ($parts.Count -gt 1) {
        $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
        $newVal = $rotated -join ", "
        $cell.Value = $newVal
    }
}
